$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume(1h) (E) cells to remain text, matching the
# inline-string storage used by the source data, then write the updated values
# scraped by the GitHub Actions job.
$cells = @{
    'D2' = '302.37'
    'E2' = '-0.64%'
    'D3' = '37.26'
    'E3' = '6.45%'
    'D4' = '4.999'
    'E4' = '-3.13%'
    'D5' = '0.07824'
    'E5' = '0.77%'
    'D6' = '2.209'
    'E6' = '-3.56%'
    'D7' = '8.031'
    'E7' = '-0.11%'
    'D8' = '4.031'
    'E8' = '0.38%'
    'D9' = '0.9144'
    'E9' = '-1.52%'
    'D10' = '0.09682'
    'E10' = '-5.15%'
    'D11' = '0.1887'
    'E11' = '3.08%'
    'D12' = '0.08669'
    'E12' = '0.08%'
    'D13' = '0.03523'
    'E13' = '1.75%'
    'D14' = '0.09957'
    'E14' = '0.99%'
    'D15' = '0.001481'
    'E15' = '-0.74%'
    'D16' = '0.005638'
    'E16' = '-3.24%'
    'D17' = '3.459'
    'E17' = '-1.35%'
    'D18' = '2.262'
    'E18' = '7.54%'
    'E19' = '1.25%'
    'D20' = '0.1300'
    'E20' = '-2.09%'
    'D21' = '4.760'
    'E21' = '2.91%'
    'D22' = '0.2295'
    'E22' = '-0.14%'
    'D23' = '0.04637'
    'E23' = '0.68%'
    'D24' = '0.001231'
    'E24' = '0.19%'
    'D25' = '0.004789'
    'E25' = '8.26%'
    'E26' = '-7.84%'
    'E27' = '38.91%'
    'D39' = '0.01767'
    'E39' = '-0.01%'
    'D40' = '0.04742'
    'E40' = '0.19%'
    'D41' = '0.008065'
    'E41' = '5.75%'
    'D42' = '0.1391'
    'E42' = '-1.07%'
    'D43' = '0.007663'
    'E43' = '7.90%'
    'D44' = '0.002211'
    'E44' = '-0.88%'
    'D45' = '0.009869'
    'E45' = '7.22%'
    'D46' = '0.00006018'
    'E46' = '1.92%'
    'D47' = '0.00000000751'
    'E47' = '-0.14%'
    'D48' = '7.926'
    'E48' = '190.09%'
    'E49' = '-0.54%'
    'E50' = '-0.14%'
    'E51' = '-0.14%'
}

foreach ($addr in $cells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cells[$addr]
}
